$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 491.7857
$ws.Range("I12").Value = 437.3846
$ws.Range("K12").Value = 437.3846
$ws.Range("M12").Value = -267.3846
$ws.Range("H32").Value = 4629.364
$ws.Range("I32").Value = 2282.1667
$ws.Range("J32").Value = 5509.5625
$ws.Range("K32").Value = 2282.1667
$ws.Range("L32").Value = 5509.5625
$ws.Range("M32").Value = -1956.1667
$ws.Range("N32").Value = -6161.5625
$ws.Range("H40").Value = 838620.9399999999
$ws.Range("I40").Value = 1432065
$ws.Range("K40").Value = 1432065
$ws.Range("M40").Value = -1431890
$ws.Range("H62").Value = 7500
$ws.Range("I62").Value = 7500
$ws.Range("K62").Value = 7500
$ws.Range("M62").Value = -6876
$ws.Range("H65").Value = 7500
$ws.Range("I65").Value = 7500
$ws.Range("K65").Value = 37500
$ws.Range("M65").Value = -34380
$ws.Range("H101").Value = 543.2143
$ws.Range("I101").Value = 539.3333
$ws.Range("K101").Value = 1617.9999
$ws.Range("M101").Value = 4.000099999999975
$ws.Range("H118").Value = 742
$ws.Range("I118").Value = 742
$ws.Range("K118").Value = 2226
$ws.Range("M118").Value = -569
$ws.Range("H129").Value = 2236.889
$ws.Range("I129").Value = 1783
$ws.Range("J129").Value = 3144.6667
$ws.Range("K129").Value = 5349
$ws.Range("L129").Value = 9434.000100000001
$ws.Range("M129").Value = -349
$ws.Range("N129").Value = -19434.0001
$ws.Range("H138").Value = 47621636
$ws.Range("I138").Value = 1632.1333
$ws.Range("J138").Value = 166671650
$ws.Range("K138").Value = 4896.3999
$ws.Range("L138").Value = 500014950
$ws.Range("M138").Value = 243.6000999999997
$ws.Range("N138").Value = -500025230

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4849.9673
$ws.Range("I32").Value = 4930
$ws.Range("K32").Value = 4930
$ws.Range("M32").Value = -4643
$ws.Range("H61").Value = 2949.8462
$ws.Range("I61").Value = 2229.2068
$ws.Range("J61").Value = 5039.7
$ws.Range("K61").Value = 2229.2068
$ws.Range("L61").Value = 5039.7
$ws.Range("M61").Value = -2017.2068
$ws.Range("N61").Value = -5463.7
$ws.Range("H74").Value = 25346.459
$ws.Range("I74").Value = 29235.244
$ws.Range("J74").Value = 2569.2856
$ws.Range("K74").Value = 29235.244
$ws.Range("L74").Value = 2569.2856
$ws.Range("M74").Value = -28361.244
$ws.Range("N74").Value = -4317.2856
$ws.Range("H77").Value = 25346.459
$ws.Range("I77").Value = 29235.244
$ws.Range("J77").Value = 2569.2856
$ws.Range("K77").Value = 146176.22
$ws.Range("L77").Value = 12846.428
$ws.Range("M77").Value = -141808.22
$ws.Range("N77").Value = -21582.428
$ws.Range("H97").Value = 1129.6154
$ws.Range("I97").Value = 1190.909
$ws.Range("J97").Value = 792.5
$ws.Range("K97").Value = 1190.909
$ws.Range("L97").Value = 792.5
$ws.Range("M97").Value = -694.9090000000001
$ws.Range("N97").Value = -1784.5
$ws.Range("H122").Value = 3409.2942
$ws.Range("I122").Value = 1452.2727
$ws.Range("K122").Value = 4356.8181
$ws.Range("M122").Value = -1906.8181
$ws.Range("H132").Value = 62268.074
$ws.Range("I132").Value = 1814.1482
$ws.Range("J132").Value = 187826.23
$ws.Range("K132").Value = 5442.444600000001
$ws.Range("L132").Value = 563478.6900000001
$ws.Range("M132").Value = -2912.444600000001
$ws.Range("N132").Value = -568538.6900000001
$ws.Range("H136").Value = 2949.8462
$ws.Range("I136").Value = 2229.2068
$ws.Range("J136").Value = 5039.7
$ws.Range("K136").Value = 6687.6204
$ws.Range("L136").Value = 15119.1
$ws.Range("M136").Value = -4137.6204
$ws.Range("N136").Value = -20219.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2111.5715
$ws.Range("I94").Value = 1704.5333
$ws.Range("K94").Value = 1704.5333
$ws.Range("M94").Value = -1253.5333
$ws.Range("H105").Value = 9911.691999999999
$ws.Range("I105").Value = 11350.182
$ws.Range("K105").Value = 11350.182
$ws.Range("M105").Value = -9603.182000000001
$ws.Range("H113").Value = 9208
$ws.Range("I113").Value = 9208
$ws.Range("K113").Value = 9208
$ws.Range("M113").Value = -7038
$ws.Range("H134").Value = 3754.7058
$ws.Range("I134").Value = 2294.6924
$ws.Range("J134").Value = 8499.75
$ws.Range("K134").Value = 6884.0772
$ws.Range("L134").Value = 25499.25
$ws.Range("M134").Value = -4349.0772
$ws.Range("N134").Value = -30569.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 116692220
$ws.Range("J4").Value = 1000000000
$ws.Range("L4").Value = 1000000000
$ws.Range("N4").Value = -1000000224
$ws.Range("H58").Value = 3324.5
$ws.Range("I58").Value = 2674.6667
$ws.Range("K58").Value = 2674.6667
$ws.Range("M58").Value = -2471.6667
$ws.Range("H122").Value = 2940.7856
$ws.Range("I122").Value = 2038.8572
$ws.Range("J122").Value = 3842.7144
$ws.Range("K122").Value = 6116.571599999999
$ws.Range("L122").Value = 11528.1432
$ws.Range("M122").Value = -3666.571599999999
$ws.Range("N122").Value = -16428.1432
$ws.Range("H134").Value = 2552.1292
$ws.Range("I134").Value = 1940.3077
$ws.Range("K134").Value = 5820.9231
$ws.Range("M134").Value = -3285.9231
$ws.Range("H136").Value = 3324.5
$ws.Range("I136").Value = 2674.6667
$ws.Range("K136").Value = 8024.000100000001
$ws.Range("M136").Value = -5474.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.666668
$ws.Range("I2").Value = 63.944443
$ws.Range("J2").Value = 38.833332
$ws.Range("K2").Value = 383.666658
$ws.Range("L2").Value = 232.999992
$ws.Range("M2").Value = -270.666658
$ws.Range("N2").Value = -458.999992
$ws.Range("H7").Value = 197.625
$ws.Range("I7").Value = 63.333332
$ws.Range("K7").Value = 189.999996
$ws.Range("M7").Value = -77.99999600000001
$ws.Range("H14").Value = 9439
$ws.Range("I14").Value = 9439
$ws.Range("K14").Value = 28317
$ws.Range("M14").Value = -28144
$ws.Range("H122").Value = 1701.0526
$ws.Range("J122").Value = 1724.7646
$ws.Range("L122").Value = 15522.8814
$ws.Range("N122").Value = -20422.8814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2533.6
$ws.Range("I97").Value = 2277.6667
$ws.Range("K97").Value = 2277.6667
$ws.Range("M97").Value = -1781.6667
$ws.Range("H122").Value = 1290.9474
$ws.Range("I122").Value = 1309.4615
$ws.Range("K122").Value = 3928.3845
$ws.Range("M122").Value = -1478.3845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3151.8823
$ws.Range("I7").Value = 1762.7142
$ws.Range("J7").Value = 4124.3
$ws.Range("K7").Value = 1762.7142
$ws.Range("L7").Value = 4124.3
$ws.Range("M7").Value = -1650.7142
$ws.Range("N7").Value = -4348.3
$ws.Range("H22").Value = 1706.375
$ws.Range("I22").Value = 1535.2
$ws.Range("J22").Value = 1991.6666
$ws.Range("K22").Value = 1535.2
$ws.Range("L22").Value = 1991.6666
$ws.Range("M22").Value = -1240.2
$ws.Range("N22").Value = -2581.6666
$ws.Range("H27").Value = 1706.375
$ws.Range("I27").Value = 1535.2
$ws.Range("J27").Value = 1991.6666
$ws.Range("K27").Value = 1535.2
$ws.Range("L27").Value = 1991.6666
$ws.Range("M27").Value = -1428.2
$ws.Range("N27").Value = -2205.6666
$ws.Range("H55").Value = 545.9048
$ws.Range("I55").Value = 196.66667
$ws.Range("K55").Value = 196.66667
$ws.Range("M55").Value = -23.66667000000001
$ws.Range("H122").Value = 6361.6
$ws.Range("I122").Value = 3943.8
$ws.Range("J122").Value = 8779.4
$ws.Range("K122").Value = 11831.4
$ws.Range("L122").Value = 26338.2
$ws.Range("M122").Value = -9381.400000000001
$ws.Range("N122").Value = -31238.2
$ws.Range("H126").Value = 3151.8823
$ws.Range("I126").Value = 1762.7142
$ws.Range("J126").Value = 4124.3
$ws.Range("K126").Value = 5288.142599999999
$ws.Range("L126").Value = 12372.9
$ws.Range("M126").Value = -2818.142599999999
$ws.Range("N126").Value = -17312.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 56112828
$ws.Range("I100").Value = 67335110
$ws.Range("J100").Value = 1397.6666
$ws.Range("K100").Value = 134670220
$ws.Range("L100").Value = 2795.3332
$ws.Range("M100").Value = -134669679
$ws.Range("N100").Value = -3877.3332
$ws.Range("H122").Value = 57479.832
$ws.Range("I122").Value = 144042.86
$ws.Range("J122").Value = 2394.2727
$ws.Range("K122").Value = 432128.58
$ws.Range("L122").Value = 7182.8181
$ws.Range("M122").Value = -429678.58
$ws.Range("N122").Value = -12082.8181
$ws.Range("H136").Value = 11439.939
$ws.Range("I136").Value = 10253.056
$ws.Range("J136").Value = 12864.2
$ws.Range("K136").Value = 30759.168
$ws.Range("L136").Value = 38592.60000000001
$ws.Range("M136").Value = -28209.168
$ws.Range("N136").Value = -43692.60000000001
